$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------
# Merge the run that was split by the stray "_GoBack" bookmark
# (", Einfüh" | bookmark | "rung in automatische Tests, Objektgleichheit: ")
# back into a single run and drop the bookmark.
$d.Content.Find.Execute(", Einführung in automatische Tests, Objektgleichheit: ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    ", Einführung in automatische Tests, Objektgleichheit: ", 2) | Out-Null

# --- Change 2 -------------------------------------------------------
# Append two new rows to the end of the topics table.
$t = $d.Tables.Item(1)

$row1 = $t.Rows.Add()
$row1.Cells.Item(1).Range.Text = "12.05."
$row1.Cells.Item(2).Range.Text = "Persistenz, Klassen zum direkten Dateizugriff (File, FileInfo, Directory, DirectoryInfo, DriveInfo), Encoding, Streaming, Objektgraphen serialisieren, XML"

$row2 = $t.Rows.Add()
$row2.Cells.Item(1).Range.Text = "19.05."
$row2.Cells.Item(2).Range.Text = "Deserialisierung, Reflection, Einführung in WPF, Visuelle Elemente (Visuals, Controls, Panels, Items Controls), Master-Detail-Ansichten, Formularansichten"

# The "_GoBack" bookmark (left behind by the last edit position) now
# belongs at the very end of the document/table.
$end = $d.Content.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($end, $end))
